{"js": "// Append five new list items to the end of the \"Quick start guide\" document:\n//   - a top-level (\"Shut everything down when you're done\") bullet\n//   - four second-level bullets underneath it, one of which (\"Type sudo\n//     shutdown now\") keeps the \"sudo\" word wrapped in spell-check\n//     <w:proofErr> markers just like the rest of the document.\n//\n// We build the new content as a WordprocessingML (\"flat OPC\") fragment and\n// insert it with Range.insertOoxml so the resulting markup (proofErr runs,\n// numPr/ilvl, etc.) matches exactly, then insert it right after the last\n// paragraph in the document body.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst insertionRange = lastParagraph.getRange(Word.RangeLocation.end);\n\nconst newContentOoxml = `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"0\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n            </w:pPr>\n            <w:r>\n              <w:t>Shut everything down when you\\u2019re done</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"1\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n            </w:pPr>\n            <w:r>\n              <w:t xml:space=\"preserve\">Type </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>sudo</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\"> shutdown now</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"1\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n            </w:pPr>\n            <w:r>\n              <w:t>Unplug battery from drone</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"1\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n            </w:pPr>\n            <w:r>\n              <w:t>Switch Vicon off</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"1\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n            </w:pPr>\n            <w:r>\n              <w:t>Remove table and test bench from flight room</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ninsertionRange.insertOoxml(newContentOoxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Append five new list items to the end of the \"Quick start guide\" document:\n#   - a top-level (\"Shut everything down when you're done\") bullet\n#   - four second-level bullets underneath it, one of which (\"Type sudo\n#     shutdown now\") keeps the \"sudo\" word wrapped in spell-check\n#     <w:proofErr> markers just like the rest of the document.\n#\n# The new content is expressed as a WordprocessingML (\"flat OPC\") fragment\n# and inserted with Range.InsertXML at a range collapsed to the very end of\n# the document, so the resulting markup (proofErr runs, numPr/ilvl, etc.)\n# matches exactly and the existing content is left untouched.\n\n$d = $word.ActiveDocument\n\n$insertionRange = $d.Content\n$insertionRange.Collapse(0)  # wdCollapseEnd\n\n$newContentOoxml = @'\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"0\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n            </w:pPr>\n            <w:r>\n              <w:t>Shut everything down when you\u2019re done</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"1\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n            </w:pPr>\n            <w:r>\n              <w:t xml:space=\"preserve\">Type </w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:t>sudo</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n            <w:r>\n              <w:t xml:space=\"preserve\"> shutdown now</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"1\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n            </w:pPr>\n            <w:r>\n              <w:t>Unplug battery from drone</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"1\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n            </w:pPr>\n            <w:r>\n              <w:t>Switch Vicon off</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListParagraph\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"1\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n            </w:pPr>\n            <w:r>\n              <w:t>Remove table and test bench from flight room</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$insertionRange.InsertXML($newContentOoxml)\n"}
